# "phone survey available after 3 weeks and some textual changes"
# Update the TapCounter app naming (now explicitly "QuantActions TapCounter")
# and fix the "Iphone" -> "iPhone" capitalisation typo in the register.xlsx
# translation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# have_smartphone row (row 4): fix "Iphone" -> "iPhone" casing
$ws.Range("B4").Value = "I have a smartphone running Android (not an iPhone)"
$ws.Range("C4").Value = "Ik heb een smartphone met Android (geen iPhone)"

# app_installed row (row 5): rename TapCounter -> 'QuantActions TapCounter'
$ws.Range("B5").Value = "I have the 'QuantActions TapCounter' app installed on my smartphone"
$ws.Range("C5").Value = "Ik heb de 'QuantActions TapCounter' app geïnstalleerd op mijn smartphone"

# button1 row (row 9): rename TapCounter -> 'QuantActions TapCounter'
$ws.Range("B9").Value = "Read more about how to install the 'QuantActions TapCounter' app"
$ws.Range("C9").Value = "Lees meer over het installeren van de 'QuantActions TapCounter' app"

# Column C now holds long Dutch strings too (e.g. row 60) -- widen/autofit it,
# same as columns A and B already are, so the text isn't clipped.
$ws.Columns.Item(3).AutoFit()
